$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Unveiling Quantum Entanglement: A Mysterious Dance of Particles" "The Journey Through the Realm of Biology: Unveiling the Secrets of Life"

# --- Author name ---
Replace-Text " Emily Carter" " Sophia Kennedy"

# --- Email line: 5 runs "emily" / "." / "carter@quantamverse" / "." / "edu"
#     become 3 runs "sophiakennedy@educonnect" / "." / "org" ---
$emailPara = $d.Paragraphs(3)
$runs = $emailPara.Range

# First run text "emily" -> "sophiakennedy@educonnect"
$full = $d.Content.Text
$start = $full.IndexOf("emily")
$len = "emily".Length
$r1 = $d.Range($start, $start + $len)
$r1.Text = "sophiakennedy@educonnect"

# Now remove the "carter@quantamverse" run entirely (it sits right after the "." run)
$full = $d.Content.Text
$start = $full.IndexOf("carter@quantamverse")
$len = "carter@quantamverse".Length
$r2 = $d.Range($start, $start + $len)
$r2.Delete()

# Now change "edu" -> "org" (the run right after the remaining ".")
$full = $d.Content.Text
$start = $full.IndexOf("edu", $start)
$len = "edu".Length
$r3 = $d.Range($start, $start + $len)
$r3.Text = "org"

# --- Body paragraph 1 (essay) ---
Replace-Text "In the realm of quantum mechanics, a realm where the laws of classical physics falter, lies a perplexing phenomenon known as quantum entanglement" "Delving into the realm of biology is akin to embarking on an extraordinary voyage through the intricate web of life, where mysteries unfold and wonders abound"

Replace-Text " This enigmatic connection between particles transcends the constraints of time and space, allowing them to share information instantaneously, regardless of the distance separating them" " Biology, the study of living organisms, unveils the captivating secrets that orchestrate the symphony of life on our planet"

Replace-Text " In this essay, we will delve into the world of quantum entanglement, unraveling its complexities and delving into its profound implications for our understanding of the universe" " It is a boundless realm where curiosity and discovery converge, beckoning us to explore the enigmatic tapestry of nature's designs"

Replace-Text "The strange and elusive nature of quantum entanglement has captivated the minds of scientists and philosophers alike since its discovery in the mid-20th century" "Biology unveils the intricate machinery that governs the very essence of life"

Replace-Text " It defies our intuition and challenges our notions of locality and causality" " From the microscopic symphony of cells to the awe-inspiring majesty of organisms, biology uncovers the profound interconnectedness that binds all living entities"

Replace-Text " As we explore this fascinating phenomenon, we will examine the experiments that have confirmed its existence, the theories that attempt to explain it, and the potential applications that it may hold for the future of computing, cryptography, and information transfer" " It delves into the hidden realms of genetics, where the blueprint of life is inscribed, and explores the marvels of evolution, where organisms adapt, thrive, and transform"

Replace-Text "Quantum entanglement has the potential to fundamentally alter our understanding of the universe" "The study of biology transcends mere knowledge acquisition; it fosters an appreciation for the diversity and unity of life"

Replace-Text " If particles can communicate instantaneously over vast distances, it raises profound questions about the nature of reality and the role of locality in the laws of physics" " Through biology, we gain insights into our own existence, unraveling the complexities of human anatomy and physiology, and fostering a profound respect for the delicate balance that sustains life on Earth"

# Collapse of: " Furthermore, ..." + "." + " Unveiling the mysteries ... transform the way we live" (3 runs)
# into a single run with new text.
$full = $d.Content.Text
$start = $full.IndexOf(" Furthermore, the ability to manipulate and harness quantum entanglement could open up new avenues for technology, revolutionizing communication, computation, and cryptography")
$endAnchor = " Unveiling the mysteries of quantum entanglement is a scientific endeavor of immense importance, with the potential to reshape our understanding of the universe and transform the way we live"
$endStart = $full.IndexOf($endAnchor)
$end = $endStart + $endAnchor.Length
$rBlock = $d.Range($start, $end)
$rBlock.Text = " It is a subject that ignites our imagination, kindles our curiosity, and inspires us to ponder the greatest mysteries of existence"

# --- Summary paragraph ---
Replace-Text "Quantum entanglement, a mysterious phenomenon in the realm of quantum mechanics, defies our classical understanding of locality and causality" "Biology unveils the intricate tapestry of life, delving into the secrets of living organisms, from the microscopic to the macroscopic"

Replace-Text " Two entangled particles, regardless of their distance apart, share information instantaneously" " It encompasses the study of genetics, evolution, and the interconnectedness of life, fostering an appreciation for diversity and unity"

# Collapse of: " Scientists have conducted..." + "enigmatic characteristics" (w/ lastRenderedPageBreak) + "." +
# " The potential applications..." + "." + " Unveiling the mysteries ... revolutionizing technology" (6 runs)
# into a single run with new text.
$full = $d.Content.Text
$start = $full.IndexOf(" Scientists have conducted experiments confirming the existence of this phenomenon, and theories have emerged to explain its ")
$endAnchor2 = " Unveiling the mysteries of quantum entanglement is a captivating and transformative scientific pursuit, holding the key to unlocking the secrets of the universe and revolutionizing technology"
$endStart2 = $full.IndexOf($endAnchor2)
$end2 = $endStart2 + $endAnchor2.Length
$rBlock2 = $d.Range($start, $end2)
$rBlock2.Text = " Biology ignites curiosity, inspires exploration, and cultivates a profound understanding of our own existence and the intricate web of life on Earth"

# --- Append a new empty paragraph at the very end of the body (before sectPr) ---
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>")

Write-Output "done"
